# The presentation's custom "Integral" design (theme1.xml, used by the
# slide master / all slides) is switched to the default "Office Theme"
# color palette (the palette that theme2.xml already carried). Only the
# <a:clrScheme> RGB values differ between the two themes - fonts and the
# format scheme (fills/lines/effects) are identical - so the edit is
# expressed as theme-color-scheme writes, which PowerPoint persists into
# the slide master's theme part (ppt/theme/theme1.xml).
#
# ThemeColorScheme slot order matches the OOXML <a:clrScheme> child order:
#   1 dk1, 2 lt1, 3 dk2, 4 lt2, 5-10 accent1-6, 11 hlink, 12 folHlink
# RGB values below are the standard Office-theme colors expressed as
# OLE COLORREF ints (0x00BBGGRR) so TextRange/Font-style ".RGB" writes
# land on the right bytes.

$p = $ppt.ActivePresentation

$officeThemeColors = @(
    0,          # dk1      000000
    16777215,   # lt1      FFFFFF
    6968388,    # dk2      44546A
    15132391,   # lt2      E7E6E6
    13998939,   # accent1  5B9BD5
    3243501,    # accent2  ED7D31
    10855845,   # accent3  A5A5A5
    49407,      # accent4  FFC000
    12874308,   # accent5  4472C4
    4697456,    # accent6  70AD47
    12673797,   # hlink    0563C1
    7491477     # folHlink 954F72
)

# The theme is shared by the whole deck (slide master level), so touching
# it through any single slide updates it for all slides/layouts.
$slide = $p.Slides.Item(1)
$colorScheme = $slide.ThemeColorScheme

for ($i = 1; $i -le $officeThemeColors.Count; $i++) {
    $colorScheme.Colors($i).RGB = $officeThemeColors[$i - 1]
}
